$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "OrderNo"
$ws.Range("B1").Value = "Product"
$ws.Range("C1").Value = "Fulfilment Store"
$ws.Range("D1").Value = "Total Price"

$ws.Range("A2").Value = 41000000006
$ws.Range("B2").Value = "Keyboard (KB123)"
$ws.Range("C2").Value = "Liverpool"
$ws.Range("D2").Value = 1000
$ws.Range("D2").NumberFormat = '"$"#,##0;[Red]\-"$"#,##0'

$ws.Columns.Item(1).ColumnWidth = 19.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

$ws.Range("F9").Select()
